$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.518.61'
$ws.Range("E2").Value = '  +0.02%  '

$ws.Range("D3").Value = '1.823.20'
$ws.Range("E3").Value = '  -0.16%  '

$ws.Range("D4").NumberFormatLocal = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("D5").NumberFormatLocal = "@"
$ws.Range("D5").Value = '315.29'
$ws.Range("E5").Value = '  -0.55%  '

$ws.Range("D7").NumberFormatLocal = "@"
$ws.Range("D7").Value = '0.5112'
$ws.Range("E7").Value = '  -5.52%  '

$ws.Range("D8").NumberFormatLocal = "@"
$ws.Range("D8").Value = '0.3949'
$ws.Range("E8").Value = '  -2.45%  '

$ws.Range("D9").NumberFormatLocal = "@"
$ws.Range("D9").Value = '0.08186'
$ws.Range("E9").Value = '  +6.80%  '

$ws.Range("D10").NumberFormatLocal = "@"
$ws.Range("D10").Value = '1.110'
$ws.Range("E10").Value = '  -0.83%  '

$ws.Range("D11").NumberFormatLocal = "@"
$ws.Range("D11").Value = '41.64'
$ws.Range("E11").Value = '  -0.50%  '

$ws.Range("D12").NumberFormatLocal = "@"
$ws.Range("D12").Value = '21.11'
$ws.Range("E12").Value = '  +0.57%  '

$ws.Range("D13").NumberFormatLocal = "@"
$ws.Range("D13").Value = '6.316'
$ws.Range("E13").Value = '  -0.15%  '

$ws.Range("D14").NumberFormatLocal = "@"
$ws.Range("D14").Value = '1.001'
$ws.Range("E14").Value = '  +0.12%  '

$ws.Range("D15").NumberFormatLocal = "@"
$ws.Range("D15").Value = '7.520'
$ws.Range("E15").Value = '  -1.55%  '

$ws.Range("D16").Value = '1.822.67'
$ws.Range("E16").Value = '  -0.11%  '

$ws.Range("E17").Value = '  +3.62%  '

$ws.Range("D18").NumberFormatLocal = "@"
$ws.Range("D18").Value = '92.63'
$ws.Range("E18").Value = '  +2.97%  '

$ws.Range("D19").NumberFormatLocal = "@"
$ws.Range("D19").Value = '0.06653'
$ws.Range("E19").Value = '  +0.77%  '

$ws.Range("E20").Value = '  +0.40%  '

$ws.Range("E21").Value = '  +0.08%  '

$ws.Range("D22").NumberFormatLocal = "@"
$ws.Range("D22").Value = '6.089'
$ws.Range("E22").Value = '  +0.28%  '

$ws.Range("D23").Value = '28.549.93'
$ws.Range("E23").Value = '  +0.09%  '

$ws.Range("D24").NumberFormatLocal = "@"
$ws.Range("D24").Value = '11.38'
$ws.Range("E24").Value = '  +1.87%  '

$ws.Range("D25").NumberFormatLocal = "@"
$ws.Range("D25").Value = '2.266'

$ws.Range("E26").Value = '  +3.34%  '

$ws.Range("D27").NumberFormatLocal = "@"
$ws.Range("D27").Value = '156.59'
$ws.Range("E27").Value = '  -0.79%  '

$ws.Range("D28").Value = '2.033.68'
$ws.Range("E28").Value = '  -0.06%  '

$ws.Range("D29").NumberFormatLocal = "@"
$ws.Range("D29").Value = '2.401'
$ws.Range("E29").Value = '  -1.91%  '

$ws.Range("E30").Value = '  +1.65%  '

$ws.Range("D31").NumberFormatLocal = "@"
$ws.Range("D31").Value = '1.112'
$ws.Range("E31").Value = '  -1.11%  '

$ws.Range("D32").NumberFormatLocal = "@"
$ws.Range("D32").Value = '0.1090'
$ws.Range("E32").Value = '  -1.38%  '

$ws.Range("D33").NumberFormatLocal = "@"
$ws.Range("D33").Value = '5.752'
$ws.Range("E33").Value = '  +1.30%  '

$ws.Range("D34").NumberFormatLocal = "@"
$ws.Range("D34").Value = '3.660'
$ws.Range("E34").Value = '  +0.55%  '

$ws.Range("D35").NumberFormatLocal = "@"
$ws.Range("D35").Value = '0.07037'
$ws.Range("E35").Value = '  -5.32%  '

$ws.Range("E36").Value = '  -0.46%  '

$ws.Range("D37").NumberFormatLocal = "@"
$ws.Range("D37").Value = '5.285'

$ws.Range("D38").NumberFormatLocal = "@"
$ws.Range("D38").Value = '0.02346'
$ws.Range("E38").Value = '  +0.00%  '

$ws.Range("D39").NumberFormatLocal = "@"
$ws.Range("D39").Value = '8.844'
$ws.Range("E39").Value = '  -0.49%  '

$ws.Range("D40").NumberFormatLocal = "@"
$ws.Range("D40").Value = '0.6313'
$ws.Range("E40").Value = '  +0.26%  '

$ws.Range("E41").Value = '  -0.61%  '

$ws.Range("D42").NumberFormatLocal = "@"
$ws.Range("D42").Value = '1.183'
$ws.Range("E42").Value = '  -0.45%  '

$ws.Range("D43").NumberFormatLocal = "@"
$ws.Range("D43").Value = '1.001'
$ws.Range("E43").Value = '  +0.12%  '

$ws.Range("D44").NumberFormatLocal = "@"
$ws.Range("D44").Value = '1.401'
$ws.Range("E44").Value = '  +0.22%  '

$ws.Range("D45").NumberFormatLocal = "@"
$ws.Range("D45").Value = '13.54'
$ws.Range("E45").Value = '  +0.64%  '

$ws.Range("D46").NumberFormatLocal = "@"
$ws.Range("D46").Value = '0.5927'
$ws.Range("E46").Value = '  +0.91%  '

$ws.Range("D47").NumberFormatLocal = "@"
$ws.Range("D47").Value = '3.731'
$ws.Range("E47").Value = '  +0.85%  '

$ws.Range("D48").NumberFormatLocal = "@"
$ws.Range("D48").Value = '125.19'
$ws.Range("E48").Value = '  -0.26%  '

$ws.Range("D49").NumberFormatLocal = "@"
$ws.Range("D49").Value = '1.989'
$ws.Range("E49").Value = '  -0.76%  '

$ws.Range("E50").Value = '  -0.86%  '

$ws.Range("D51").NumberFormatLocal = "@"
$ws.Range("D51").Value = '0.06905'
$ws.Range("E51").Value = '  +0.23%  '
